# This patch removes internal rounding of floating point values used to
# compute the simulation outputs. As a result, the numeric values stored
# in the reference output sheet change from integers to full-precision
# floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17359.6299122635
$ws.Range("C2").Value = 36891.45214766602
$ws.Range("D2").Value = 71785.74973061308
$ws.Range("E2").Value = 108831.3431461096

# Row 3
$ws.Range("B3").Value = 179130.1737333099
$ws.Range("C3").Value = 332107.6342700651
$ws.Range("D3").Value = 402421.5175193835
$ws.Range("E3").Value = 448190.6277148937

# Row 4
$ws.Range("B4").Value = 19233.0375699881
$ws.Range("C4").Value = 36301.32689215025
$ws.Range("D4").Value = 57786.38454577468
$ws.Range("E4").Value = 75011.34798940629

# Row 6
$ws.Range("B6").Value = 105189.0445423489
$ws.Range("C6").Value = 133734.4124316712
$ws.Range("D6").Value = 122851.7835609764
$ws.Range("E6").Value = 100780.2597267853

# Row 7
$ws.Range("B7").Value = 10868.29414806273
$ws.Range("C7").Value = 21744.58556810208
$ws.Range("D7").Value = 23933.45715564557
$ws.Range("E7").Value = 26067.72939541312

# Row 9
$ws.Range("B9").Value = 809341.8732435276
$ws.Range("C9").Value = 1298027.669546911
$ws.Range("D9").Value = 1757763.466731737
$ws.Range("E9").Value = 2167541.737071982

# Row 12
$ws.Range("B12").Value = 788885.2271656395
$ws.Range("C12").Value = 867423.9748857918
$ws.Range("D12").Value = 725474.0045054841
$ws.Range("E12").Value = 547087.2432607213
